$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the date/time number format used in column B (rows 2:35) down to
# the newly added rows (36:46) before writing their values.
$ws.Range("B36:B46").NumberFormat = $ws.Range("B35").NumberFormat

$ws.Cells.Item(2, 1).Value = 5631
$ws.Cells.Item(2, 2).Value = 45992.95833333334
$ws.Cells.Item(3, 1).Value = 5543
$ws.Cells.Item(3, 2).Value = 45992.96875
$ws.Cells.Item(4, 1).Value = 5452
$ws.Cells.Item(4, 2).Value = 45992.97916666666
$ws.Cells.Item(5, 1).Value = 5331
$ws.Cells.Item(5, 2).Value = 45992.98958333334
$ws.Cells.Item(6, 1).Value = 5269
$ws.Cells.Item(6, 2).Value = 45993
$ws.Cells.Item(7, 1).Value = 5272
$ws.Cells.Item(7, 2).Value = 45993.01041666666
$ws.Cells.Item(8, 1).Value = 5223
$ws.Cells.Item(8, 2).Value = 45993.02083333334
$ws.Cells.Item(9, 1).Value = 5145
$ws.Cells.Item(9, 2).Value = 45993.03125
$ws.Cells.Item(10, 1).Value = 5158
$ws.Cells.Item(10, 2).Value = 45993.04166666666
$ws.Cells.Item(11, 1).Value = 5131
$ws.Cells.Item(11, 2).Value = 45993.05208333334
$ws.Cells.Item(12, 1).Value = 5183
$ws.Cells.Item(12, 2).Value = 45993.0625
$ws.Cells.Item(13, 1).Value = 5110
$ws.Cells.Item(13, 2).Value = 45993.07291666666
$ws.Cells.Item(14, 1).Value = 5159
$ws.Cells.Item(14, 2).Value = 45993.08333333334
$ws.Cells.Item(15, 1).Value = 5177
$ws.Cells.Item(15, 2).Value = 45993.09375
$ws.Cells.Item(16, 1).Value = 5151
$ws.Cells.Item(16, 2).Value = 45993.10416666666
$ws.Cells.Item(17, 1).Value = 5097
$ws.Cells.Item(17, 2).Value = 45993.11458333334
$ws.Cells.Item(18, 1).Value = 5176
$ws.Cells.Item(18, 2).Value = 45993.125
$ws.Cells.Item(19, 1).Value = 5219
$ws.Cells.Item(19, 2).Value = 45993.13541666666
$ws.Cells.Item(20, 1).Value = 5222
$ws.Cells.Item(20, 2).Value = 45993.14583333334
$ws.Cells.Item(21, 1).Value = 5266
$ws.Cells.Item(21, 2).Value = 45993.15625
$ws.Cells.Item(22, 1).Value = 5375
$ws.Cells.Item(22, 2).Value = 45993.16666666666
$ws.Cells.Item(23, 1).Value = 5439
$ws.Cells.Item(23, 2).Value = 45993.17708333334
$ws.Cells.Item(24, 1).Value = 5488
$ws.Cells.Item(24, 2).Value = 45993.1875
$ws.Cells.Item(25, 1).Value = 5591
$ws.Cells.Item(25, 2).Value = 45993.19791666666
$ws.Cells.Item(26, 1).Value = 5947
$ws.Cells.Item(26, 2).Value = 45993.20833333334
$ws.Cells.Item(27, 1).Value = 6101
$ws.Cells.Item(27, 2).Value = 45993.21875
$ws.Cells.Item(28, 1).Value = 6345
$ws.Cells.Item(28, 2).Value = 45993.22916666666
$ws.Cells.Item(29, 1).Value = 6524
$ws.Cells.Item(29, 2).Value = 45993.23958333334
$ws.Cells.Item(30, 1).Value = 6914
$ws.Cells.Item(30, 2).Value = 45993.25
$ws.Cells.Item(31, 1).Value = 7122
$ws.Cells.Item(31, 2).Value = 45993.26041666666
$ws.Cells.Item(32, 1).Value = 7291
$ws.Cells.Item(32, 2).Value = 45993.27083333334
$ws.Cells.Item(33, 1).Value = 7368
$ws.Cells.Item(33, 2).Value = 45993.28125
$ws.Cells.Item(34, 1).Value = 7561
$ws.Cells.Item(34, 2).Value = 45993.29166666666
$ws.Cells.Item(35, 1).Value = 7674
$ws.Cells.Item(35, 2).Value = 45993.30208333334
$ws.Cells.Item(36, 1).Value = 7792
$ws.Cells.Item(36, 2).Value = 45993.3125
$ws.Cells.Item(37, 1).Value = 7781
$ws.Cells.Item(37, 2).Value = 45993.32291666666
$ws.Cells.Item(38, 1).Value = 7842
$ws.Cells.Item(38, 2).Value = 45993.33333333334
$ws.Cells.Item(39, 1).Value = 7874
$ws.Cells.Item(39, 2).Value = 45993.34375
$ws.Cells.Item(40, 1).Value = 7979
$ws.Cells.Item(40, 2).Value = 45993.35416666666
$ws.Cells.Item(41, 1).Value = 7905
$ws.Cells.Item(41, 2).Value = 45993.36458333334
$ws.Cells.Item(42, 1).Value = 7869
$ws.Cells.Item(42, 2).Value = 45993.375
$ws.Cells.Item(43, 1).Value = 7834
$ws.Cells.Item(43, 2).Value = 45993.38541666666
$ws.Cells.Item(44, 1).Value = 7821
$ws.Cells.Item(44, 2).Value = 45993.39583333334
$ws.Cells.Item(45, 1).Value = 7755
$ws.Cells.Item(45, 2).Value = 45993.40625
$ws.Cells.Item(46, 1).Value = 7741
$ws.Cells.Item(46, 2).Value = 45993.41666666666
